# thesis metadata.xlsx — "semenov et al 2007 and tidying up metadata extraction file"
#
# - removes the min_time / max_time rows (and their definitions), replacing
#   them with a single "period_flux" row
# - adds a new "size" field (with its 0,1,2,3 code and "very small, small,
#   medium, large" meaning) at the end of the table
# - moves the reviewer comment from C23 to C22 (it now sits on the exp_age
#   row instead of the row that used to hold that content)
# - updates the saved sheet view (scroll position / selection)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- move the existing comment before we start shuffling rows around -------
# (C23 currently holds "pre-juvenile, juvenile, adult?" — the note belongs on
# the exp_age definition, which after the edit lives one row higher, at C22)
$existingComment = $ws.Range("C23").Comment
if ($existingComment) {
    $commentText = $existingComment.Text()
    $existingComment.Delete()
    $ws.Range("C22").AddComment($commentText)
}

# --- remove min_time (row 11) + max_time (row 12) ---------------------------
$ws.Rows("11:12").Delete()

# --- insert the replacement "period_flux" row at the top of that block -----
$ws.Rows("11:11").Insert()
$ws.Range("A11").Value = "period_flux"

# --- append the new "size" row after exp_age (now row 22) ------------------
$ws.Rows("23:23").Insert()
$ws.Range("A23").Value = "size"
$ws.Range("B23").Value = "0,1,2,3"
$ws.Range("C23").Value = "very small, small, medium, large "
$ws.Rows("23:23").RowHeight = 17

# --- restore the saved scroll position / selection --------------------------
$ws.Application.Goto($ws.Range("A16"), $true)
$ws.Range("G20").Select()
